# Doing Updates for Financials
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AZPN")

# Balance Sheet updates (column D only)
$ws.Range("D43").Value = 366700   # Net Receivables
$ws.Range("D46").Value = 476000   # Total Current Assets
$ws.Range("D47").Value = 340600   # Long Term Investments (was "NA")
$ws.Range("D49").Value = 112200   # Goodwill
$ws.Range("D54").Value = 950300   # Total Assets
$ws.Range("D59").Value = 56400    # Other Current Liabilities
$ws.Range("D60").Value = 230600   # Total Current Liabilities
$ws.Range("D62").Value = 243500   # Other Liabilities
$ws.Range("D66").Value = 474100   # Total Liabilities
$ws.Range("D72").Value = 1065500  # Retained Earnings
$ws.Range("D76").Value = 476200   # Total Stockholder Equity

# Cash Flow Statement updates - Capital Expenditures (row 91), columns D:J
$ws.Range("D91").Value = -300
$ws.Range("E91").Value = -2700
$ws.Range("F91").Value = -3500
$ws.Range("G91").Value = -7600
$ws.Range("H91").Value = -4000
$ws.Range("I91").Value = -4500
$ws.Range("J91").Value = -4200
